# Mejoras en redacción de puntos en la metodología
# Slide 9 ("Metodología") - reposition the two connector shapes and
# tweak a few bullet / label texts.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# ---- Reposition the "flowChartMerge" connector shapes -------------------
# (Left/Top are expressed in points = EMU/12700; a tiny epsilon is added
# so that the float32 points round-trip lands on the exact target EMU
# instead of truncating one unit short.)
$emu = 12700
$bump = 0.00002

# Shape 6: "Diagrama de flujo: combinar 20" (id 21)
$shMerge1 = $s.Shapes.Item(6)
$shMerge1.Left = 2374480 / $emu + $bump
$shMerge1.Top  = 2862056 / $emu + $bump

# Shape 7: "Diagrama de flujo: combinar 21" (id 22)
$shMerge2 = $s.Shapes.Item(7)
$shMerge2.Left = 2374480 / $emu + $bump
$shMerge2.Top  = 4819451 / $emu + $bump

# ---- Text tweaks ----------------------------------------------------------
# (Edit the existing Run's Text directly rather than the Paragraph's, so a
# single run is kept instead of being split into several runs.)

# Shape 8: "Rectángulo: esquinas redondeadas 4" (id 5)
$shBox1 = $s.Shapes.Item(8)
$shBox1.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Factores Críticos para ZOFRI S.A."

# Shape 9: "Rectángulo: esquinas redondeadas 22" (id 23)
$shBox2 = $s.Shapes.Item(9)
$shBox2.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Diseño de la Estructura"
$shBox2.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Construcción de los Contenidos"

# Shape 10: "Rectángulo: esquinas redondeadas 23" (id 24)
$shBox3 = $s.Shapes.Item(10)
$shBox3.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Definición de las actividades"
